$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend header numbers D1:G1, matching B1/C1 style ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1:G1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5

# --- Row 2 ---
$ws.Range("B2").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': '001-SEM (18450)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['T', 'R'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"
$ws.Range("C2").Value = "{'name': 'ARABL-UH 1120 - Elementary Arabic 2', 'title': '003-SEM (18353)', 'id': 'ARABLUH1120160332', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 11.833333333333, 'end_date': 13.083333333333, 'inscturct_mode': 'P', 'instructor': 'Muhamed Al-Khalil', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("D2").Value = "{'name': 'ARABL-UH 2120 - Intermediate Arabic 2', 'title': '001-SEM (18354)', 'id': 'ARABLUH2120204522', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 10.416666666667, 'end_date': 11.666666666667, 'inscturct_mode': 'P', 'instructor': 'Khulood Kittaneh, Omima El Araby', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("E2").Value = "{'name': 'ARTH-UH 2128 - Money and Art in the Global Renaissance', 'title': '001-SEM (22629)', 'id': 'ARTHUH2128232572', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 13.25, 'end_date': 14.5, 'inscturct_mode': 'P', 'instructor': 'Mahnaz Yousefzadeh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("F2").Value = "{'name': 'AW-UH 1118 - Archaeology, Arabia and the Bible', 'title': '001-SEM (24761)', 'id': 'AWUH1118236369', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 14.666666666667, 'end_date': 15.916666666667, 'inscturct_mode': 'P', 'instructor': 'William Zimmerle', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("G2").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': 'REC1-RCT (18468)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['U'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"

# --- Row 3 ---
$ws.Range("B3").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': '001-SEM (18450)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['T', 'R'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"
$ws.Range("C3").Value = "{'name': 'ARABL-UH 1120 - Elementary Arabic 2', 'title': '003-SEM (18353)', 'id': 'ARABLUH1120160332', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 11.833333333333, 'end_date': 13.083333333333, 'inscturct_mode': 'P', 'instructor': 'Muhamed Al-Khalil', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("D3").Value = "{'name': 'ARABL-UH 2120 - Intermediate Arabic 2', 'title': '001-SEM (18354)', 'id': 'ARABLUH2120204522', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 10.416666666667, 'end_date': 11.666666666667, 'inscturct_mode': 'P', 'instructor': 'Khulood Kittaneh, Omima El Araby', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("E3").Value = "{'name': 'ARTH-UH 2128 - Money and Art in the Global Renaissance', 'title': '001-SEM (22629)', 'id': 'ARTHUH2128232572', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 13.25, 'end_date': 14.5, 'inscturct_mode': 'P', 'instructor': 'Mahnaz Yousefzadeh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("F3").Value = "{'name': 'AW-UH 1118 - Archaeology, Arabia and the Bible', 'title': '001-SEM (24761)', 'id': 'AWUH1118236369', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 14.666666666667, 'end_date': 15.916666666667, 'inscturct_mode': 'P', 'instructor': 'William Zimmerle', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("G3").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': 'REC2-RCT (18469)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['U'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Staff', 'status': 'Wait List', 'waitlist_count': '4', 'session': 'A71'}"

# --- Row 4 ---
$ws.Range("B4").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': '001-SEM (18450)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['T', 'R'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"
$ws.Range("C4").Value = "{'name': 'ARABL-UH 1120 - Elementary Arabic 2', 'title': '002-SEM (18021)', 'id': 'ARABLUH1120160332', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 10.416666666667, 'end_date': 11.666666666667, 'inscturct_mode': 'P', 'instructor': 'Muhamed Al-Khalil', 'status': 'Wait List', 'waitlist_count': '0', 'session': 'AD'}"
$ws.Range("D4").Value = "{'name': 'ARABL-UH 2120 - Intermediate Arabic 2', 'title': '002-SEM (18786)', 'id': 'ARABLUH2120204522', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 11.833333333333, 'end_date': 13.083333333333, 'inscturct_mode': 'P', 'instructor': 'Khulood Kittaneh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("E4").Value = "{'name': 'ARTH-UH 2128 - Money and Art in the Global Renaissance', 'title': '001-SEM (22629)', 'id': 'ARTHUH2128232572', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 13.25, 'end_date': 14.5, 'inscturct_mode': 'P', 'instructor': 'Mahnaz Yousefzadeh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("F4").Value = "{'name': 'AW-UH 1118 - Archaeology, Arabia and the Bible', 'title': '001-SEM (24761)', 'id': 'AWUH1118236369', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 14.666666666667, 'end_date': 15.916666666667, 'inscturct_mode': 'P', 'instructor': 'William Zimmerle', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("G4").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': 'REC1-RCT (18468)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['U'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"

# --- Row 5 ---
$ws.Range("B5").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': '001-SEM (18450)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['T', 'R'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Dania Zantout', 'status': 'Open', 'waitlist_count': 0, 'session': 'A71'}"
$ws.Range("C5").Value = "{'name': 'ARABL-UH 1120 - Elementary Arabic 2', 'title': '002-SEM (18021)', 'id': 'ARABLUH1120160332', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 10.416666666667, 'end_date': 11.666666666667, 'inscturct_mode': 'P', 'instructor': 'Muhamed Al-Khalil', 'status': 'Wait List', 'waitlist_count': '0', 'session': 'AD'}"
$ws.Range("D5").Value = "{'name': 'ARABL-UH 2120 - Intermediate Arabic 2', 'title': '002-SEM (18786)', 'id': 'ARABLUH2120204522', 'term': '1224', 'campus': 'AD', 'days': ['M', 'T', 'W', 'R'], 'start_date': 11.833333333333, 'end_date': 13.083333333333, 'inscturct_mode': 'P', 'instructor': 'Khulood Kittaneh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("E5").Value = "{'name': 'ARTH-UH 2128 - Money and Art in the Global Renaissance', 'title': '001-SEM (22629)', 'id': 'ARTHUH2128232572', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 13.25, 'end_date': 14.5, 'inscturct_mode': 'P', 'instructor': 'Mahnaz Yousefzadeh', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("F5").Value = "{'name': 'AW-UH 1118 - Archaeology, Arabia and the Bible', 'title': '001-SEM (24761)', 'id': 'AWUH1118236369', 'term': '1224', 'campus': 'AD', 'days': ['M', 'W'], 'start_date': 14.666666666667, 'end_date': 15.916666666667, 'inscturct_mode': 'P', 'instructor': 'William Zimmerle', 'status': 'Open', 'waitlist_count': 0, 'session': 'AD'}"
$ws.Range("G5").Value = "{'name': 'MATH-UH 1000A - Mathematics for Statistics and Calculus Part I', 'title': 'REC2-RCT (18469)', 'id': 'MATHUH1000A234160', 'term': '1224', 'campus': 'AD', 'days': ['U'], 'start_date': 9.0, 'end_date': 10.25, 'inscturct_mode': 'P', 'instructor': 'Staff', 'status': 'Wait List', 'waitlist_count': '4', 'session': 'A71'}"

# --- Remove now-unused rows 6:9 entirely ---
$ws.Range("A6:A9").EntireRow.Delete() | Out-Null

Write-Host "done"